$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates (Apob -> Lsr, Target cluster = ECs)
$ws.Range("M2").Value = 0.6235436666666666
$ws.Range("N2").Value = 1.870631
$ws.Range("O2").Value = 0.3577087286482158
$ws.Range("P2").Value = 0.3577087286482157
$ws.Range("Q2").Value = 0.009050944169555555
$ws.Range("R2").Value = 0.081458497526
$ws.Range("S2").Value = 0.3577087286482158
$ws.Range("T2").Value = 0.3577087286482157

# Row 3 updates (Apob -> Lsr, Target cluster = FAPs)
$ws.Range("O3").Value = 0.5045751194047746
$ws.Range("P3").Value = 0.5045751194047745
$ws.Range("S3").Value = 0.5045751194047746
$ws.Range("T3").Value = 0.5045751194047745

# Row 4 updates (Apob -> Lsr, Target cluster = MuSCs)
$ws.Range("N4").Value = 0.7201839999999999
$ws.Range("O4").Value = 0.1377161519470097
$ws.Range("P4").Value = 0.1377161519470096
$ws.Range("S4").Value = 0.1377161519470097
$ws.Range("T4").Value = 0.1377161519470096
